$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to text format so numeric-looking
# strings (e.g. "1.001", "29.275.15") are preserved verbatim as text
# instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.275.15"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "1.911.98"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "321.13"
$ws.Range("E5").Value = "  -2.97%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.4724"
$ws.Range("E7").Value = "  +2.85%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "47.68"
$ws.Range("D10").Value = "0.08042"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("D11").Value = "1.003"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "22.48"
$ws.Range("E12").Value = "  +3.85%  "
$ws.Range("D13").Value = "1.966.76"
$ws.Range("E13").Value = "  +3.18%  "
$ws.Range("D14").Value = "5.893"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").Value = "7.129"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "89.70"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "0.06641"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "0.00001030"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "17.69"
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "29.294.27"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("D23").Value = "5.519"
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("D25").Value = "2.196"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").Value = "2.131.03"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "155.10"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").Value = "19.76"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").Value = "5.994"
$ws.Range("E29").Value = "  +10.90%  "
$ws.Range("D30").Value = "2.105"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").Value = "117.56"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "1.073"
$ws.Range("E32").Value = "  +10.09%  "
$ws.Range("D33").Value = "0.09514"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("D34").Value = "1.423"
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("D35").Value = "3.540"
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("E36").Value = "  +2.24%  "
$ws.Range("D37").Value = "0.06079"
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("D38").Value = "0.02246"
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").Value = "8.241"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").Value = "1.175"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("D41").Value = "0.5853"
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("D42").Value = "2.500"
$ws.Range("E42").Value = "  +11.54%  "
$ws.Range("D43").Value = "0.1835"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").Value = "10.10"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("D45").Value = "0.07889"
$ws.Range("E45").Value = "  +4.80%  "
$ws.Range("D46").Value = "1.270"
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("D47").Value = "12.12"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("D48").Value = "0.5524"
$ws.Range("D49").Value = "1.922"
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("D50").Value = "113.00"
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("D51").Value = "44.22"
$ws.Range("E51").Value = "  -2.24%  "

# Restore the default cell style (no explicit style index) to match
# the original workbook formatting, since setting NumberFormat above
# would otherwise leave a new style applied to these cells.
$dataRange.Style = "Normal"

